$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "NO OF HOURS LATE" (column F) computation for the 02-13-2015 row (row 8)
# was wrong - it was left blank instead of showing the 0.75 hrs late value.
# Fixing the number, and also dropping the tardy-row highlight (style changes
# from the highlighted style used for rows 5-8 to the plain style used for
# the rest of the sheet - same format already applied to row 9) now that the
# correct (non-zero) late value is reflected in the data itself.

# Copy the (already-correct) formatting from row 9 onto row 8.
$ws.Range("A9:J9").Copy()
$ws.Range("A8:J8").PasteSpecial(-4122)

# Fix the number of hours late for row 8.
$ws.Range("F8").Value = 0.75
